$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.706699999999991
$ws.Range("A8").Value = -21.1449
$ws.Range("A10").Value = -20.52049999999997
$ws.Range("A12").Value = -22.48030000000003
$ws.Range("B13").Value = 6.2655
$ws.Range("A18").Value = -22.36790000000002
$ws.Range("D20").Value = -8.160399999999999
$ws.Range("A25").Value = -22.27090000000003
